# Insert a new data row above the current row 10 (Femacal de La Calera /
# Frambuesa weekly price sheet). This pushes the existing rows 10-16 down
# to rows 11-17 and adds one new record at row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row before row 10; everything currently at row 10
# and below shifts down by one row (old row 10 -> 11, ..., old row 16 -> 17).
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the new observation.
$ws.Cells.Item(10, 1).Value  = 3
$ws.Cells.Item(10, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(10, 3).Value  = "Coquimbo"
$ws.Cells.Item(10, 4).Value  = 44960
$ws.Cells.Item(10, 5).Value  = 5
$ws.Cells.Item(10, 6).Value  = "Fruta"
$ws.Cells.Item(10, 7).Value  = 100101
$ws.Cells.Item(10, 8).Value  = "Berries"
$ws.Cells.Item(10, 9).Value  = 100101004
$ws.Cells.Item(10, 10).Value = "Frambuesa"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 40
$ws.Cells.Item(10, 14).Value = 7000
$ws.Cells.Item(10, 15).Value = 7000
$ws.Cells.Item(10, 16).Value = 7000
$ws.Cells.Item(10, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(10, 19).Value = 3500
$ws.Cells.Item(10, 20).Value = 2
